$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-09-15"

# Update the September column header label
$ws.Range("A10").Value = "September (through 09-15)"

# Update September row (row 10) values for columns C-H (2016-2021)
$ws.Range("C10").Value = 26
$ws.Range("D10").Value = 37
$ws.Range("E10").Value = 29
$ws.Range("F10").Value = 32
$ws.Range("G10").Value = 57
$ws.Range("H10").Value = 76

# Update Total row (row 11) values for columns C-H (2016-2021)
$ws.Range("C11").Value = 407
$ws.Range("D11").Value = 588
$ws.Range("E11").Value = 519
$ws.Range("F11").Value = 381
$ws.Range("G11").Value = 841
$ws.Range("H11").Value = 1146
